# Apply the "test case.xlsx" edit:
#   - Sheet2 gets a new value in A1 ("www") plus a new selection (B5).
#   - Sheet1 ("parameters") gains six new rows (A3:A8) below the existing
#     "poolname" row, and the active selection moves to C6.
#
# Shared-string interning order matters: the target sharedStrings.xml has
# "www" (from Sheet2!A1) landing at index 11, immediately followed by the
# Sheet1 values w/e/e/r/t/ty at indices 12-16. So Sheet2!A1 is written
# first, then Sheet1's new cells top-to-bottom.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "parameters"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet2"

# Sheet2: new value, interned first so it becomes shared-string index 11.
$ws2.Range("A1").Value = "www"

# Sheet1: six new rows under the existing data.
$ws1.Range("A3").Value = "w"
$ws1.Range("A4").Value = "e"
$ws1.Range("A5").Value = "e"
$ws1.Range("A6").Value = "r"
$ws1.Range("A7").Value = "t"
$ws1.Range("A8").Value = "ty"

# Match the saved selections from the diff.
$ws1.Range("C6").Select() | Out-Null
$ws2.Range("B5").Select() | Out-Null

# Leave "parameters" as the active/visible sheet (tabSelected stays there).
$ws1.Activate() | Out-Null
